$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add new column K ("Tempo empresa") to the worksheet.
# ---------------------------------------------------------------------------

# --- Header cell K1: reuse J1's cell format (style index 1) then set text ---
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)
$ws.Range("K1").Value = "Tempo empresa"

# --- New numeric data for rows 166-186 (the "Tempo empresa" observations) ---
$values = @{
    166 = 0
    167 = 3.5
    168 = 1.6
    169 = 0.5
    170 = 1
    171 = 0.16
    172 = 1.1599999999999999
    173 = 2.5
    174 = 0.75
    175 = 0.25
    176 = 3
    177 = 12
    178 = 7
    179 = 2.2999999999999998
    180 = 0
    181 = 1.3
    182 = 4
    183 = 30
    184 = 9
    185 = 0.3
    186 = 0.25
}

$ws.Range("J1").Copy()
foreach ($row in 166..186) {
    $cell = $ws.Cells.Item($row, 11)
    $cell.PasteSpecial(-4122)
    $cell.Value = $values[$row]
}

# --- Column width for the new column (closest achievable match to the
#     author's 15.5546875 character width, reusing column J's width) ---
$ws.Columns.Item(11).ColumnWidth = $ws.Columns.Item(10).ColumnWidth

# Clear marching-ants / clipboard state left over from the Copy() calls
$excel.CutCopyMode = 0

# --- Mirror the author's final selection state ---
$ws.Range("L1:L1048576").Select()
